$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 12.9
$ws.Range("D3").Value = -7.761
$ws.Range("E5").Value = 13.098
$ws.Range("D14").Value = -8.111000000000001
$ws.Range("D16").Value = -7.947999999999999
$ws.Range("E16").Value = 13.041
$ws.Range("D21").Value = -7.934
$ws.Range("D23").Value = -7.826000000000001
$ws.Range("D25").Value = -8.306999999999999
